$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'61.201.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.95%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.012.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.71%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'571.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.44%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'128.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.54%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'3.014.84"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.51%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  -1.79%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = "'  -3.55%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'5.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.84%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "'  -4.61%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.98%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'33.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.12%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "'  +0.07%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'3.509.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.72%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'61.307.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.74%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'3.013.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.70%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'6.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.60%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'439.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.32%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'  -4.78%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "'  -5.03%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = "'  -4.93%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'79.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.36%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'12.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.41%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.15%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.12%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'2.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -6.07%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = "'  -2.00%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'7.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.69%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'6.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -8.27%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'25.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.66%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'0.0947"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -8.03%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = "'  -3.84%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'0.961"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -6.45%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  -5.11%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'50.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.81%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.0₃0684"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.54%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.0367"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.98%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'7.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.25%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'  -2.87%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'373.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.73%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'2.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -8.95%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'2.656.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.84%  "
$ws.Range("E44").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'122.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.49%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'  -4.87%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'33.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.54%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  -7.26%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'  -2.91%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'23.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -6.72%  "
$ws.Range("E51").Style = "Normal"
